# Remove 4 resolved/duplicate claim rows from the INCO sheet.
# Rows identified by their "Caso" value in column A:
#   -270 (row 37)  - SALTA SUR /ALT/ 917
#   -178 (row 19)  - USPALLATA /ALT/ 3504
#   -161 (row 15)  - URQUIZA 1692
#   -145 (row 14)  - HUMBERTO PRIMO /ALT/ 2849
#
# Deleting from the bottom up keeps the remaining row numbers stable while
# we work, so no re-lookup of shifted positions is needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(37).Delete()
$ws.Rows(19).Delete()
$ws.Rows(15).Delete()
$ws.Rows(14).Delete()
